$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.516.41'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.580.95'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.71'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.498'
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.22'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('E9').Value = '  -0.76%  '
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.805.63'
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.595.44'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  -2.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.502.95'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.99'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '215.15'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0691'
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.30'
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.75'
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.99'
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.06'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('E26').Value = '  +2.55%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.16'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0472'
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.362.77'
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.969'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0168'
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.531'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.822'
$ws.Range('E40').Value = '  +0.97%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.972'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +3.64%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.76'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.26'
$ws.Range('E46').Value = '  -2.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.717.16'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.40'
$ws.Range('E48').Value = '  -1.51%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0957'
$ws.Range('E50').Value = '  -1.42%  '
$ws.Range('E51').Value = '  -0.61%  '
